$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D16").Value = "2016-03-08 06:58:05"
$wsZhCn.Range("G16").Value = "2016-03-08 06:58:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D16").Value = "2016-03-08 06:58:15"
$wsDeDe.Range("G16").Value = "2016-03-08 06:59:04"
